# Natmi following Dr Hou advice
#
# The ligand-receptor summary table originally only modelled two sending
# clusters (FAPs, sCs) against four target clusters. Per Dr Hou's advice,
# "ECs" is added as a third sending cluster, which also recomputes every
# specificity/weight metric for all sending x target cluster combinations
# (now 3 x 4 = 12 data rows instead of 2 x 4 = 8).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data set for A2:T13 - columns are:
# Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# Ligand-expressing cells, Ligand detection rate, Ligand average expression
# value, Ligand total expression value, Ligand derived specificity of
# average expression value, Ligand derived specificity of total expression
# value, Receptor-expressing cells, Receptor detection rate, Receptor
# average expression value, Receptor total expression value, Receptor
# derived specificity of average expression value, Receptor derived
# specificity of total expression value, Edge average expression weight,
# Edge total expression weight, Edge average expression derived
# specificity, Edge total expression derived specificity
$data = @(
    @("ECs", "Ctf1", "Lifr", "ECs", 1, 0.3333333333333333, 0.1923443333333333, 0.577033, 0.06576397757907125, 0.06576397757907126, 3, 1, 31.44605333333334, 94.33816, 0.273208187120734, 0.273208187120734, 6.048470164364445, 54.43623147928, 0.01796725709222665, 0.01796725709222666),
    @("ECs", "Ctf1", "Lifr", "FAPs", 1, 0.3333333333333333, 0.1923443333333333, 0.577033, 0.06576397757907125, 0.06576397757907126, 3, 1, 57.80064033333333, 173.401921, 0.5021809252974907, 0.5021809252974908, 11.11762563115478, 100.058630680393, 0.03302541511190143, 0.03302541511190145),
    @("ECs", "Ctf1", "Lifr", "M2", 1, 0.3333333333333333, 0.1923443333333333, 0.577033, 0.06576397757907125, 0.06576397757907126, 3, 1, 12.04190666666667, 36.12572, 0.1046219522368387, 0.1046219522368387, 2.316192509862223, 20.84573258876, 0.006880355721182121, 0.006880355721182123),
    @("ECs", "Ctf1", "Lifr", "sCs", 1, 0.3333333333333333, 0.1923443333333333, 0.577033, 0.06576397757907125, 0.06576397757907126, 3, 1, 13.81063466666667, 41.431904, 0.1199889353449366, 0.1199889353449366, 2.656397317870222, 23.907575860832, 0.007890949653761041, 0.007890949653761043),
    @("FAPs", "Ctf1", "Lifr", "ECs", 3, 1, 1.496366666666667, 4.4891, 0.5116190438852002, 0.5116190438852003, 3, 1, 31.44605333333334, 94.33816, 0.273208187120734, 0.273208187120734, 47.05482600622222, 423.493434056, 0.1397785114763188, 0.1397785114763188),
    @("FAPs", "Ctf1", "Lifr", "FAPs", 3, 1, 1.496366666666667, 4.4891, 0.5116190438852002, 0.5116190438852003, 3, 1, 57.80064033333333, 173.401921, 0.5021809252974907, 0.5021809252974908, 86.49095150678887, 778.4185635610999, 0.2569253248580873, 0.2569253248580874),
    @("FAPs", "Ctf1", "Lifr", "M2", 3, 1, 1.496366666666667, 4.4891, 0.5116190438852002, 0.5116190438852003, 3, 1, 12.04190666666667, 36.12572, 0.1046219522368387, 0.1046219522368387, 18.01910773911111, 162.171969652, 0.05352658317281448, 0.05352658317281449),
    @("FAPs", "Ctf1", "Lifr", "sCs", 3, 1, 1.496366666666667, 4.4891, 0.5116190438852002, 0.5116190438852003, 3, 1, 13.81063466666667, 41.431904, 0.1199889353449366, 0.1199889353449366, 20.66577336071111, 185.9919602464, 0.06138862437797958, 0.06138862437797959),
    @("sCs", "Ctf1", "Lifr", "ECs", 3, 1, 1.236056333333333, 3.708169, 0.4226169785357285, 0.4226169785357286, 3, 1, 31.44605333333334, 94.33816, 0.273208187120734, 0.273208187120734, 38.86909338100444, 349.82184042904, 0.1154624185521885, 0.1154624185521886),
    @("sCs", "Ctf1", "Lifr", "FAPs", 3, 1, 1.236056333333333, 3.708169, 0.4226169785357285, 0.4226169785357286, 3, 1, 57.80064033333333, 173.401921, 0.5021809252974907, 0.5021809252974908, 71.44484755473876, 643.003627992649, 0.2122301853275019, 0.212230185327502),
    @("sCs", "Ctf1", "Lifr", "M2", 3, 1, 1.236056333333333, 3.708169, 0.4226169785357285, 0.4226169785357286, 3, 1, 12.04190666666667, 36.12572, 0.1046219522368387, 0.1046219522368387, 14.88447500074222, 133.96027500668, 0.04421501334284206, 0.04421501334284207),
    @("sCs", "Ctf1", "Lifr", "sCs", 3, 1, 1.236056333333333, 3.708169, 0.4226169785357285, 0.4226169785357286, 3, 1, 13.81063466666667, 41.431904, 0.1199889353449366, 0.1199889353449366, 17.07072244708622, 153.636502023776, 0.050709361313196, 0.05070936131319601),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    for ($c = 1; $c -le 20; $c++) {
        $ws.Cells.Item($r, $c).Value = $data[$i][$c - 1]
    }
}
